# Apply the nightly cryptos-list refresh (see commit message: "Updated
# cryptos list ... with GitHub Actions"). For every changed row, update the
# Price (column D) and Volume(1h) (column E) cells with their new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values look like plain decimal numbers (e.g. "1.000",
# "7.420", "0.000008051"). Excel would normally auto-convert such strings to
# numeric values when assigned, which silently strips trailing zeros or
# switches to scientific notation. Force those specific cells to Text format
# first so the values are stored verbatim as text, matching the sheet's
# existing convention of keeping Price/Volume figures as text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.253.20"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.924.13"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "248.36"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "0.7152"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.3185"
$ws.Range("E8").Value = "  -4.13%  "
$ws.Range("D9").Value = "27.71"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "0.07052"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").Value = "0.7904"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "0.07971"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "1.927.40"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "5.371"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "94.68"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "14.66"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "30.246.31"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "256.74"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "0.000008051"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").Value = "5.761"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "2.180.97"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "6.827"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "9.522"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "165.42"
$ws.Range("E26").Value = "  +3.83%  "
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "2.255"
$ws.Range("E28").Value = "  -8.48%  "
$ws.Range("D29").Value = "0.1257"
$ws.Range("E29").Value = "  -6.27%  "
$ws.Range("D30").Value = "1.354"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").Value = "4.384"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "4.112"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").Value = "0.05126"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "1.265"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").Value = "0.7427"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "2.763"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").Value = "0.01952"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("D40").Value = "77.28"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Value = "6.347"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("D42").Value = "0.4498"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "1.986"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "0.8439"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "100.42"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("D47").Value = "9.716"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "7.420"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("D49").Value = "36.49"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "0.06171"
$ws.Range("E50").Value = "  +3.60%  "
$ws.Range("D51").Value = "0.4195"
$ws.Range("E51").Value = "  +2.50%  "
